$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Commands")
$ws.Activate()

# Insert a new row right after row 95 (the "Artisan Command" section header row)
# for the new "alarm(n,<bool>)" command, pushing every following row down by one.
$ws.Rows.Item(96).Insert()

$ws.Range("B96").Value = "alarm(n,<bool>)"
$ws.Range("C96").Value = "enables/disables alarm number n"

# Style the description text in C96: the trailing "n" is italic, the rest is not.
$head = $ws.Range("C96").Characters(1, 30)
$head.Font.Name = "Calibri"
$head.Font.Size = 11
$head.Font.Color = 0

$tail = $ws.Range("C96").Characters(31, 1)
$tail.Font.Name = "Calibri"
$tail.Font.Size = 11
$tail.Font.Color = 0
$tail.Font.Italic = $true

$ws.Range("C98").Select()
